$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.874.97'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '2.313.04'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '108.72'
$ws.Range("E5").Value = '  +10.26%  '
$ws.Range("D6").Value = '271.87'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("D10").Value = '48.16'
$ws.Range("E10").Value = '  +5.10%  '
$ws.Range("D11").Value = '0.0943'
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("E12").Value = '  +5.23%  '
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("D14").Value = '15.82'
$ws.Range("E14").Value = '  +1.71%  '
$ws.Range("D15").Value = '2.652.96'
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").Value = '0.866'
$ws.Range("E16").Value = '  -1.37%  '
$ws.Range("D17").Value = '2.303.99'
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").Value = '43.835.03'
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("E19").Value = '  +1.92%  '
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("D21").Value = '72.43'
$ws.Range("E21").Value = '  -1.73%  '
$ws.Range("D22").Value = '2.51'
$ws.Range("E22").Value = '  +8.61%  '
$ws.Range("D23").Value = '234.34'
$ws.Range("E23").Value = '  -2.56%  '
$ws.Range("D24").Value = '2.96'
$ws.Range("E24").Value = '  +16.29%  '
$ws.Range("D25").Value = '9.34'
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '11.41'
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").Value = '41.03'
$ws.Range("E28").Value = '  +6.53%  '
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("D30").Value = '2.29'
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("D31").Value = '177.47'
$ws.Range("E31").Value = '  +1.04%  '
$ws.Range("D32").Value = '21.95'
$ws.Range("E32").Value = '  -2.21%  '
$ws.Range("D33").Value = '0.0918'
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("E35").Value = '  +7.67%  '
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("D38").Value = '0.0358'
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("D39").Value = '3.88'
$ws.Range("E39").Value = '  +14.22%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.239'
$ws.Range("E40").Value = '  -2.56%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").Value = '2.36'
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("E42").Value = '  -1.62%  '
$ws.Range("D43").Value = '67.34'
$ws.Range("D44").Value = '12.26'
$ws.Range("E44").Value = '  -1.68%  '
$ws.Range("D45").Value = '5.60'
$ws.Range("E45").Value = '  +4.59%  '
$ws.Range("D46").Value = '8.81'
$ws.Range("E46").Value = '  -3.60%  '
$ws.Range("E47").Value = '  -1.78%  '
$ws.Range("E48").Value = '  +1.42%  '
$ws.Range("D49").Value = '99.65'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("E50").Value = '  +4.40%  '
$ws.Range("D51").Value = '2.547.49'
$ws.Range("E51").Value = '  -0.19%  '
